# Updates LR-pair metrics for Hras-Cav1 with recalculated TPM-based values.
# Cell values below are taken from the new (post-TPM-update) NATMI output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.535856000000001
$ws.Range("H2").Value = 25.607568
$ws.Range("I2").Value = 0.36987004643386
$ws.Range("J2").Value = 0.36987004643386
$ws.Range("M2").Value = 788.1599833333333
$ws.Range("N2").Value = 2364.47995
$ws.Range("O2").Value = 0.8397951873720987
$ws.Range("P2").Value = 0.8397951873720988
$ws.Range("Q2").Value = 6727.620122695734
$ws.Range("R2").Value = 60548.5811042616
$ws.Range("S2").Value = 0.3106150849482502
$ws.Range("T2").Value = 0.3106150849482503

# Row 3
$ws.Range("G3").Value = 8.535856000000001
$ws.Range("H3").Value = 25.607568
$ws.Range("I3").Value = 0.36987004643386
$ws.Range("J3").Value = 0.36987004643386
$ws.Range("O3").Value = 0.01890163353898316
$ws.Range("P3").Value = 0.01890163353898317
$ws.Range("Q3").Value = 151.4214561607627
$ws.Range("R3").Value = 1362.793105446864
$ws.Range("S3").Value = 0.006991148074739508
$ws.Range("T3").Value = 0.006991148074739509

# Row 4
$ws.Range("G4").Value = 8.535856000000001
$ws.Range("H4").Value = 25.607568
$ws.Range("I4").Value = 0.36987004643386
$ws.Range("J4").Value = 0.36987004643386
$ws.Range("M4").Value = 131.4690986666667
$ws.Range("N4").Value = 394.407296
$ws.Range("O4").Value = 0.1400821136357036
$ws.Range("P4").Value = 0.1400821136357036
$ws.Range("Q4").Value = 1122.201294668459
$ws.Range("R4").Value = 10099.81165201613
$ws.Range("S4").Value = 0.05181217787499094
$ws.Range("T4").Value = 0.05181217787499094

# Row 5
$ws.Range("G5").Value = 8.535856000000001
$ws.Range("H5").Value = 25.607568
$ws.Range("I5").Value = 0.36987004643386
$ws.Range("J5").Value = 0.36987004643386
$ws.Range("M5").Value = 1.145987666666667
$ws.Range("N5").Value = 3.437963
$ws.Range("O5").Value = 0.001221065453214498
$ws.Range("P5").Value = 0.001221065453214498
$ws.Range("Q5").Value = 9.781985700442668
$ws.Range("R5").Value = 88.03787130398401
$ws.Range("S5").Value = 0.0004516355358792285
$ws.Range("T5").Value = 0.0004516355358792285

# Row 6
$ws.Range("I6").Value = 0.2236685002562326
$ws.Range("J6").Value = 0.2236685002562326
$ws.Range("M6").Value = 788.1599833333333
$ws.Range("N6").Value = 2364.47995
$ws.Range("O6").Value = 0.8397951873720987
$ws.Range("P6").Value = 0.8397951873720988
$ws.Range("Q6").Value = 4068.338914289689
$ws.Range("R6").Value = 36615.05022860719
$ws.Range("S6").Value = 0.1878357300819192
$ws.Range("T6").Value = 0.1878357300819192

# Row 7
$ws.Range("I7").Value = 0.2236685002562326
$ws.Range("J7").Value = 0.2236685002562326
$ws.Range("O7").Value = 0.01890163353898316
$ws.Range("P7").Value = 0.01890163353898317
$ws.Range("S7").Value = 0.004227700026057271
$ws.Range("T7").Value = 0.004227700026057272

# Row 8
$ws.Range("I8").Value = 0.2236685002562326
$ws.Range("J8").Value = 0.2236685002562326
$ws.Range("M8").Value = 131.4690986666667
$ws.Range("N8").Value = 394.407296
$ws.Range("O8").Value = 0.1400821136357036
$ws.Range("P8").Value = 0.1400821136357036
$ws.Range("Q8").Value = 678.6196475874417
$ws.Range("R8").Value = 6107.576828286975
$ws.Range("S8").Value = 0.03133195626962098
$ws.Range("T8").Value = 0.03133195626962098

# Row 9
$ws.Range("I9").Value = 0.2236685002562326
$ws.Range("J9").Value = 0.2236685002562326
$ws.Range("M9").Value = 1.145987666666667
$ws.Range("N9").Value = 3.437963
$ws.Range("O9").Value = 0.001221065453214498
$ws.Range("P9").Value = 0.001221065453214498
$ws.Range("Q9").Value = 5.915380529569777
$ws.Range("R9").Value = 53.238424766128
$ws.Range("S9").Value = 0.0002731138786351836
$ws.Range("T9").Value = 0.0002731138786351837

# Row 10
$ws.Range("G10").Value = 5.335438
$ws.Range("H10").Value = 16.006314
$ws.Range("I10").Value = 0.2311916579666972
$ws.Range("J10").Value = 0.2311916579666973
$ws.Range("M10").Value = 788.1599833333333
$ws.Range("N10").Value = 2364.47995
$ws.Range("O10").Value = 0.8397951873720987
$ws.Range("P10").Value = 0.8397951873720988
$ws.Range("Q10").Value = 4205.178725156034
$ws.Range("R10").Value = 37846.6085264043
$ws.Range("S10").Value = 0.1941536417210087
$ws.Range("T10").Value = 0.1941536417210087

# Row 11
$ws.Range("G11").Value = 5.335438
$ws.Range("H11").Value = 16.006314
$ws.Range("I11").Value = 0.2311916579666972
$ws.Range("J11").Value = 0.2311916579666973
$ws.Range("O11").Value = 0.01890163353898316
$ws.Range("P11").Value = 0.01890163353898317
$ws.Range("Q11").Value = 94.64777653412466
$ws.Range("R11").Value = 851.8299888071219
$ws.Range("S11").Value = 0.004369899996156449
$ws.Range("T11").Value = 0.004369899996156451

# Row 12
$ws.Range("G12").Value = 5.335438
$ws.Range("H12").Value = 16.006314
$ws.Range("I12").Value = 0.2311916579666972
$ws.Range("J12").Value = 0.2311916579666973
$ws.Range("M12").Value = 131.4690986666667
$ws.Range("N12").Value = 394.407296
$ws.Range("O12").Value = 0.1400821136357036
$ws.Range("P12").Value = 0.1400821136357036
$ws.Range("Q12").Value = 701.4452248518827
$ws.Range("R12").Value = 6313.007023666943
$ws.Range("S12").Value = 0.0323858161029176
$ws.Range("T12").Value = 0.03238581610291761

# Row 13
$ws.Range("G13").Value = 5.335438
$ws.Range("H13").Value = 16.006314
$ws.Range("I13").Value = 0.2311916579666972
$ws.Range("J13").Value = 0.2311916579666973
$ws.Range("M13").Value = 1.145987666666667
$ws.Range("N13").Value = 3.437963
$ws.Range("O13").Value = 0.001221065453214498
$ws.Range("P13").Value = 0.001221065453214498
$ws.Range("Q13").Value = 6.114346144264666
$ws.Range("R13").Value = 55.029115298382
$ws.Range("S13").Value = 0.0002823001466145163
$ws.Range("T13").Value = 0.0002823001466145163

# Row 14
$ws.Range("G14").Value = 4.044874
$ws.Range("H14").Value = 12.134622
$ws.Range("I14").Value = 0.1752697953432102
$ws.Range("J14").Value = 0.1752697953432102
$ws.Range("M14").Value = 788.1599833333333
$ws.Range("N14").Value = 2364.47995
$ws.Range("O14").Value = 0.8397951873720987
$ws.Range("P14").Value = 0.8397951873720988
$ws.Range("Q14").Value = 3188.007824425433
$ws.Range("R14").Value = 28692.0704198289
$ws.Range("S14").Value = 0.1471907306209206
$ws.Range("T14").Value = 0.1471907306209206

# Row 15
$ws.Range("G15").Value = 4.044874
$ws.Range("H15").Value = 12.134622
$ws.Range("I15").Value = 0.1752697953432102
$ws.Range("J15").Value = 0.1752697953432102
$ws.Range("O15").Value = 0.01890163353898316
$ws.Range("P15").Value = 0.01890163353898317
$ws.Range("Q15").Value = 71.75387109000067
$ws.Range("R15").Value = 645.784839810006
$ws.Range("S15").Value = 0.003312885442029937
$ws.Range("T15").Value = 0.003312885442029938

# Row 16
$ws.Range("G16").Value = 4.044874
$ws.Range("H16").Value = 12.134622
$ws.Range("I16").Value = 0.1752697953432102
$ws.Range("J16").Value = 0.1752697953432102
$ws.Range("M16").Value = 131.4690986666667
$ws.Range("N16").Value = 394.407296
$ws.Range("O16").Value = 0.1400821136357036
$ws.Range("P16").Value = 0.1400821136357036
$ws.Range("Q16").Value = 531.7759390002346
$ws.Range("R16").Value = 4785.983451002112
$ws.Range("S16").Value = 0.02455216338817409
$ws.Range("T16").Value = 0.02455216338817409

# Row 17
$ws.Range("G17").Value = 4.044874
$ws.Range("H17").Value = 12.134622
$ws.Range("I17").Value = 0.1752697953432102
$ws.Range("J17").Value = 0.1752697953432102
$ws.Range("M17").Value = 1.145987666666667
$ws.Range("N17").Value = 3.437963
$ws.Range("O17").Value = 0.001221065453214498
$ws.Range("P17").Value = 0.001221065453214498
$ws.Range("Q17").Value = 4.635375717220667
$ws.Range("R17").Value = 41.718381454986
$ws.Range("S17").Value = 0.0002140158920855692
$ws.Range("T17").Value = 0.0002140158920855692

